# Remove the "tumble_dryer" and "washing_machine" rows from the
# Info Technology behavior-scenario list (annual electricity consumption
# reduced to 1200+ -> these two appliance rows are dropped from the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 = tumble_dryer, Row 24 = washing_machine (1-indexed sheet rows).
# Deleting row 23 twice removes both, shifting everything below up by two
# rows and renumbering the index column (A) automatically stays as typed
# values already in column A for the remaining rows below, so fix those too.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(23).Delete()

# Re-number column A (the ID_Technology index) sequentially for the rows
# that shifted up, since those are literal values rather than a formula.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

$ws.Application.CutCopyMode = $false

# Reset the view so it isn't scrolled down to where the deleted rows used
# to be, and land the selection back near the top of the list.
$ws.Application.Goto($ws.Range("A13"), $true) | Out-Null
$ws.Range("A13").Select() | Out-Null
